# Actualización 10 de Mayo
# Populate the "Rescatables" sheet with the list of rescatable students (rows 2-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A: NC (clave numérica)
$nc = @(19330051920441,19330051920321,19330051920375,19330051920377,19330051920382,
        18330051920116,18330051920120,18330051920121,18330051920125,17330051920466,
        18330051920455,19330051920307,19330051920312,19330051920337,19330051920373,
        19330051920374,19330051920381)

# Column B: Paterno
$paterno = @("GARCIA","DE JESUS","HERNANDEZ","HERNANDEZ","MAZAHUA","ESPIRITU","GIL",
             "GONZALEZ","LIBONATTI","MARROQUIN","PAZ","APALE","CARDENAS","MONTIEL",
             "GARCIA","HERNANDEZ","LOPEZ")

# Column C: Materno
$materno = @("ANTONIO","DE LA CRUZ","ANTONIO","FLORES","IXMATLAHUA","BUSTAMANTE","MARTINEZ",
             "REYES","FIGUEROA","ELIAS","MONTERROSAS","ZEPAHUA","AMADOR","ROJAS","DIAZ",
             "AGUILA","GARCIA")

# Column D: Nombres
$nombres = @("ABRAHAM","IGNACIO","MARIA GUADALUPE","PERLA","LUCERO","JOSE ANTONIO","ARELY",
             "ANGEL RODRIGO","EBERTH JUVIEL","JORGE ANTONIO","DANIEL","ISRAEL","KEVIN HONAM",
             "VICTOR YAHIR","GISELA","JESUS","MARIAM ABRIL")

# Column E: Nombre_Largo (materia)
$materia = @("FÍSICA I","FÍSICA I","FÍSICA I","FÍSICA I","FÍSICA I",
             "TEMAS DE FÍSICA","TEMAS DE FÍSICA","TEMAS DE FÍSICA","TEMAS DE FÍSICA","TEMAS DE FÍSICA","TEMAS DE FÍSICA",
             "FÍSICA I","FÍSICA I","FÍSICA I","FÍSICA I","FÍSICA I","FÍSICA I")

# Column F: Grupo
$grupo = @("4APV","4APV","4ARHV","4ARHV","4ARHV",
           "6APV","6APV","6APV","6APV","6APV","6APV",
           "4APV","4APV","4APV","4ARHV","4ARHV","4ARHV")

# Column G: Reprobadas
$reprobadas = @(2,2,2,2,2,2,2,2,2,2,2,1,1,1,1,1,1)

$count = $nc.Count

for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 1).Value = $nc[$i] }
for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 2).Value = $paterno[$i] }
for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 3).Value = $materno[$i] }
for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 4).Value = $nombres[$i] }
for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 5).Value = $materia[$i] }
for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 6).Value = $grupo[$i] }
for ($i = 0; $i -lt $count; $i++) { $ws.Cells.Item($i + 2, 7).Value = $reprobadas[$i] }
